$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 224081
$ws.Range("E2").Value = 2013
$ws.Range("F2").Value = 2013
$ws.Range("G2").Value = 373
$ws.Range("H2").Value = 316
$ws.Range("I2").Value = 321
$ws.Range("J2").Value = -5
$ws.Range("K2").Value = 85047
$ws.Range("L2").Value = 59865
$ws.Range("M2").Value = 25182
$ws.Range("N2").Value = 25256
$ws.Range("O2").Value = -74
$ws.Range("P2").Value = 6487
$ws.Range("Q2").Value = -2344
$ws.Range("R2").Value = 3382
$ws.Range("S2").Value = -1166
$ws.Range("T2").Value = 1550
$ws.Range("U2").Value = -3893
$ws.Range("V2").Value = 27394
$ws.Range("W2").Value = 0.9
$ws.Range("X2").Value = 0.14
$ws.Range("Y2").Value = 1.26
$ws.Range("Z2").Value = 0.36
$ws.Range("AA2").Value = 237.72
$ws.Range("AB2").Value = 288.91
$ws.Range("AC2").Value = 129
$ws.Range("AD2").Value = 69.22
$ws.Range("AE2").Value = 10171
$ws.Range("AF2").Value = 0.88
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 1.12
$ws.Range("AI2").Value = 77.44
$ws.Range("AJ2").Value = 248187647

# Row 3
$ws.Range("D3").Value = 195234
$ws.Range("E3").Value = 1911
$ws.Range("F3").Value = 1916
$ws.Range("G3").Value = 1243
$ws.Range("H3").Value = 730
$ws.Range("I3").Value = 735
$ws.Range("J3").Value = -5
$ws.Range("K3").Value = 82032
$ws.Range("L3").Value = 56838
$ws.Range("M3").Value = 25194
$ws.Range("N3").Value = 25278
$ws.Range("O3").Value = -84
$ws.Range("P3").Value = 6487
$ws.Range("Q3").Value = 1192
$ws.Range("R3").Value = 258
$ws.Range("S3").Value = -3065
$ws.Range("T3").Value = 1244
$ws.Range("U3").Value = -52
$ws.Range("V3").Value = 24730
$ws.Range("W3").Value = 0.98
$ws.Range("X3").Value = 0.37
$ws.Range("Y3").Value = 2.91
$ws.Range("Z3").Value = 0.87
$ws.Range("AA3").Value = 225.6
$ws.Range("AB3").Value = 296.44
$ws.Range("AC3").Value = 296
$ws.Range("AD3").Value = 18.54
$ws.Range("AE3").Value = 10180
$ws.Range("AF3").Value = 0.54
$ws.Range("AG3").Value = 100
$ws.Range("AH3").Value = 1.82
$ws.Range("AI3").Value = 33.78
$ws.Range("AJ3").Value = 248187647

# Row 4
$ws.Range("D4").Value = 129047
$ws.Range("E4").Value = 1554
$ws.Range("F4").Value = 1673
$ws.Range("G4").Value = 707
$ws.Range("H4").Value = -816
$ws.Range("I4").Value = -817
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 86203
$ws.Range("L4").Value = 62148
$ws.Range("M4").Value = 24055
$ws.Range("N4").Value = 24280
$ws.Range("O4").Value = -225
$ws.Range("P4").Value = 6487
$ws.Range("Q4").Value = 146
$ws.Range("R4").Value = -5633
$ws.Range("S4").Value = 1177
$ws.Range("T4").Value = 1042
$ws.Range("U4").Value = -896
$ws.Range("V4").Value = 27448
$ws.Range("W4").Value = 1.21
$ws.Range("X4").Value = -0.63
$ws.Range("Y4").Value = -3.3
$ws.Range("Z4").Value = -0.97
$ws.Range("AA4").Value = 258.35
$ws.Range("AB4").Value = 280.51
$ws.Range("AC4").Value = -329
$ws.Range("AD4").Value = -21.04
$ws.Range("AE4").Value = 9778
$ws.Range("AF4").Value = 0.71
$ws.Range("AG4").Value = 100
$ws.Range("AH4").Value = 1.45
$ws.Range("AI4").Value = -30.4
$ws.Range("AJ4").Value = 248187647

# Row 5
$ws.Range("D5").Value = 152023
$ws.Range("E5").Value = 1428
$ws.Range("F5").Value = 1428
$ws.Range("G5").Value = 628
$ws.Range("H5").Value = 346
$ws.Range("I5").Value = 334
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 72274
$ws.Range("L5").Value = 48789
$ws.Range("M5").Value = 23485
$ws.Range("N5").Value = 23702
$ws.Range("O5").Value = -217
$ws.Range("P5").Value = 6487
$ws.Range("Q5").Value = -6346
$ws.Range("R5").Value = 8780
$ws.Range("S5").Value = -6761
$ws.Range("T5").Value = 1567
$ws.Range("U5").Value = -7914
$ws.Range("V5").Value = 21129
$ws.Range("W5").Value = 0.9399999999999999
$ws.Range("X5").Value = 0.23
$ws.Range("Y5").Value = 1.39
$ws.Range("Z5").Value = 0.44
$ws.Range("AA5").Value = 207.74
$ws.Range("AB5").Value = 281.63
$ws.Range("AC5").Value = 135
$ws.Range("AD5").Value = 49.42
$ws.Range("AE5").Value = 9849
$ws.Range("AF5").Value = 0.68
$ws.Range("AG5").Value = 120
$ws.Range("AH5").Value = 1.8
$ws.Range("AI5").Value = 86.44
$ws.Range("AJ5").Value = 248187647

# Row 6
$ws.Range("D6").Value = 139865
$ws.Range("E6").Value = 1379
$ws.Range("F6").Value = 1379
$ws.Range("G6").Value = 601
$ws.Range("H6").Value = 77
$ws.Range("I6").Value = 73
$ws.Range("K6").Value = 77690
$ws.Range("L6").Value = 54572
$ws.Range("M6").Value = 23119
$ws.Range("N6").Value = 23366
$ws.Range("P6").Value = 6487
$ws.Range("Q6").Value = -1990
$ws.Range("R6").Value = -826
$ws.Range("S6").Value = 5038
$ws.Range("T6").Value = 1271
$ws.Range("U6").Value = -3261
$ws.Range("V6").Value = 26300
$ws.Range("W6").Value = 0.99
$ws.Range("X6").Value = 0.06
$ws.Range("Y6").Value = 0.31
$ws.Range("Z6").Value = 0.1
$ws.Range("AA6").Value = 236.05
$ws.Range("AB6").Value = 276.51
$ws.Range("AC6").Value = 29
$ws.Range("AD6").Value = 178.04
$ws.Range("AE6").Value = 9709
$ws.Range("AF6").Value = 0.54
$ws.Range("AG6").Value = 120
$ws.Range("AH6").Value = 2.31
$ws.Range("AI6").Value = 398.25
$ws.Range("AJ6").Value = 248187647

# Row 7
$ws.Range("D7").Value = 142183
$ws.Range("E7").Value = 2482
$ws.Range("G7").Value = 1131
$ws.Range("H7").Value = 718
$ws.Range("I7").Value = 646
$ws.Range("K7").Value = 97243
$ws.Range("L7").Value = 72627
$ws.Range("M7").Value = 24620
$ws.Range("N7").Value = 23730
$ws.Range("P7").Value = 6490
$ws.Range("Q7").Value = 10930
$ws.Range("R7").Value = -2430
$ws.Range("S7").Value = -603
$ws.Range("T7").Value = 1780
$ws.Range("U7").Value = 8260
$ws.Range("W7").Value = 1.75
$ws.Range("X7").Value = 0.51
$ws.Range("Y7").Value = 2.74
$ws.Range("Z7").Value = 0.82
$ws.Range("AA7").Value = 294.99
$ws.Range("AC7").Value = 260
$ws.Range("AD7").Value = 19.12
$ws.Range("AE7").Value = 9861
$ws.Range("AF7").Value = 0.5
$ws.Range("AG7").Value = 80
$ws.Range("AH7").Value = 1.61
$ws.Range("AI7").Value = 30.73

# Row 8
$ws.Range("D8").Value = 146594
$ws.Range("E8").Value = 2978
$ws.Range("G8").Value = 1702
$ws.Range("H8").Value = 1251
$ws.Range("I8").Value = 1091
$ws.Range("K8").Value = 99453
$ws.Range("L8").Value = 73807
$ws.Range("M8").Value = 25647
$ws.Range("N8").Value = 24610
$ws.Range("P8").Value = 6490
$ws.Range("Q8").Value = 10057
$ws.Range("R8").Value = -2420
$ws.Range("S8").Value = -1353
$ws.Range("T8").Value = 3800
$ws.Range("U8").Value = 13190
$ws.Range("W8").Value = 2.03
$ws.Range("X8").Value = 0.85
$ws.Range("Y8").Value = 4.52
$ws.Range("Z8").Value = 1.27
$ws.Range("AA8").Value = 287.78
$ws.Range("AC8").Value = 440
$ws.Range("AD8").Value = 11.32
$ws.Range("AE8").Value = 10227
$ws.Range("AF8").Value = 0.49
$ws.Range("AG8").Value = 80
$ws.Range("AH8").Value = 1.61
$ws.Range("AI8").Value = 18.19

# Row 9
$ws.Range("D9").Value = 153203
$ws.Range("E9").Value = 3313
$ws.Range("G9").Value = 1980
$ws.Range("H9").Value = 1443
$ws.Range("I9").Value = 1267
$ws.Range("K9").Value = 101340
$ws.Range("L9").Value = 74453
$ws.Range("M9").Value = 26887
$ws.Range("N9").Value = 25670
$ws.Range("P9").Value = 6490
$ws.Range("Q9").Value = 7553
$ws.Range("R9").Value = -2323
$ws.Range("S9").Value = -2463
$ws.Range("T9").Value = 3800
$ws.Range("U9").Value = 8300
$ws.Range("W9").Value = 2.16
$ws.Range("X9").Value = 0.9399999999999999
$ws.Range("Y9").Value = 5.04
$ws.Range("Z9").Value = 1.44
$ws.Range("AA9").Value = 276.92
$ws.Range("AC9").Value = 510
$ws.Range("AD9").Value = 9.75
$ws.Range("AE9").Value = 10668
$ws.Range("AF9").Value = 0.47
$ws.Range("AG9").Value = 80
$ws.Range("AH9").Value = 1.61
$ws.Range("AI9").Value = 15.68
